# Comments and Deaths Model Filter >10
# Shifts the dataset forward by 6 days (new day-0 = 2020-03-15 / serial 43905),
# updates the "predicted" column with refreshed model output, clears the
# "deaths" column for days 6+ (filtered out, per commit message: ">10"), and
# removes the six trailing rows that fell off the 18-day window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: date serial, days, deaths (blank = cleared), predicted
$data = @(
    @(43905, 0, 11,   11),
    @(43906, 1, 17,   15),
    @(43907, 2, 24,   22),
    @(43908, 3, 28,   31),
    @(43909, 4, 44,   44),
    @(43910, 5, 67,   63),
    @(43911, 6, $null, 89),
    @(43912, 7, $null, 126),
    @(43913, 8, $null, 177),
    @(43914, 9, $null, 250),
    @(43915, 10, $null, 353),
    @(43916, 11, $null, 499),
    @(43917, 12, $null, 703),
    @(43918, 13, $null, 992),
    @(43919, 14, $null, 1400),
    @(43920, 15, $null, 1975),
    @(43921, 16, $null, 2787),
    @(43922, 17, $null, 3931)
)

# Remove the six rows that are no longer part of the 18-day window
# (old rows 20-25, i.e. beyond what is now the last row, 19)
$ws.Range("A20:D25").EntireRow.Delete()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]

    if ($null -eq $vals[2]) {
        $ws.Cells.Item($row, 3).ClearContents()
    } else {
        $ws.Cells.Item($row, 3).Value = $vals[2]
    }

    $ws.Cells.Item($row, 4).Value = $vals[3]
}
